# Apply updated dSF (column F) values as described in the commit:
# "repull data, push all data, mean calculation"
#
# The workbook stores weekly records starting at row 2 (header in row 1).
# Column F ("dSF") previously mirrored column E ("dS0") with all zeros
# (except a couple of rows). After repulling the source data, column F
# was recalculated/pushed with its own (possibly negative) values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    3  = -4
    4  = -3
    5  = 3
    7  = 1
    8  = -4
    9  = -2
    11 = -5
    12 = 1
    15 = -2
    16 = 1
    18 = -4
    19 = 2
    20 = -5
    21 = 3
    22 = 1
    24 = 2
    25 = 9
    26 = -1
    28 = 6
    29 = -2
    30 = 7
    31 = -2
    32 = 1
    33 = 1
    34 = -6
    35 = -1
    36 = 1
    37 = 5
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
